$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 123, shifting existing rows 123-160 down to 124-161
$ws.Rows.Item(123).Insert()

# Populate the new row 123 with the new price-report entry
$ws.Cells.Item(123, 1).Value = 11
$ws.Cells.Item(123, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(123, 3).Value = "Bíobío"
$ws.Cells.Item(123, 4).Value = 44726
$ws.Cells.Item(123, 5).Value = 8
$ws.Cells.Item(123, 6).Value = "Fruta"
$ws.Cells.Item(123, 7).Value = 100109
$ws.Cells.Item(123, 8).Value = "Uva"
$ws.Cells.Item(123, 9).Value = 100109001
$ws.Cells.Item(123, 10).Value = "Uva"
$ws.Cells.Item(123, 11).Value = "Superior Seedless"
$ws.Cells.Item(123, 12).Value = "Primera"
$ws.Cells.Item(123, 13).Value = 100
$ws.Cells.Item(123, 14).Value = 11000
$ws.Cells.Item(123, 15).Value = 12000
$ws.Cells.Item(123, 16).Value = 11500
$ws.Cells.Item(123, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(123, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(123, 19).Value = 1150
$ws.Cells.Item(123, 20).Value = 10
